$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.237.81"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.27%  "

$ws.Range("D3").Value = "'1.929.02"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.18%  "

$ws.Range("E4").Value = "  +0.20%  "

$ws.Range("D5").Value = "'248.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.08%  "

$ws.Range("D6").Value = "'0.7118"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.13%  "

$ws.Range("E7").Value = "  +0.24%  "

$ws.Range("D8").Value = "'0.3210"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.29%  "

$ws.Range("D9").Value = "'27.28"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.19%  "

$ws.Range("D10").Value = "'0.07097"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.05%  "

$ws.Range("D11").Value = "'0.7918"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.95%  "

$ws.Range("D12").Value = "'0.08028"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.32%  "

$ws.Range("D13").Value = "'1.929.78"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.11%  "

$ws.Range("D14").Value = "'5.367"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.95%  "

$ws.Range("D15").Value = "'94.80"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.12%  "

$ws.Range("D16").Value = "'14.66"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.80%  "

$ws.Range("D17").Value = "'30.253.51"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.23%  "

$ws.Range("D18").Value = "'255.17"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.73%  "

$ws.Range("D19").Value = "'0.000008028"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.26%  "

$ws.Range("D20").Value = "'5.748"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.61%  "

$ws.Range("D21").Value = "'2.183.10"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.26%  "

$ws.Range("E22").Value = "  +0.21%  "

$ws.Range("D23").Value = "'1.002"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.23%  "

$ws.Range("D24").Value = "'6.818"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.88%  "

$ws.Range("D25").Value = "'9.550"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.37%  "

$ws.Range("D26").Value = "'166.16"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.36%  "

$ws.Range("D27").Value = "'19.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.40%  "

$ws.Range("D28").Value = "'2.274"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.67%  "

$ws.Range("D29").Value = "'0.1275"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.64%  "

$ws.Range("D30").Value = "'1.357"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.53%  "

$ws.Range("E31").Value = "  -1.65%  "

$ws.Range("D32").Value = "'4.390"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.29%  "

$ws.Range("D33").Value = "'4.132"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.33%  "

$ws.Range("D34").Value = "'0.05167"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.61%  "

$ws.Range("D35").Value = "'1.264"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.54%  "

$ws.Range("D36").Value = "'0.7449"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.74%  "

$ws.Range("D37").Value = "'2.770"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.67%  "

$ws.Range("D38").Value = "'0.01953"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.65%  "

$ws.Range("D39").Value = "'2.804"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.71%  "

$ws.Range("D40").Value = "'77.69"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.17%  "

$ws.Range("E41").Value = "  -4.09%  "

$ws.Range("D42").Value = "'0.4479"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.46%  "

$ws.Range("D43").Value = "'1.973"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.99%  "

$ws.Range("D44").Value = "'0.8456"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.34%  "

$ws.Range("D45").Value = "'1.001"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.16%  "

$ws.Range("D46").Value = "'100.55"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.72%  "

$ws.Range("D47").Value = "'9.673"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.19%  "

$ws.Range("D48").Value = "'7.426"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.87%  "

$ws.Range("D49").Value = "'36.35"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.39%  "

$ws.Range("D50").Value = "'0.06124"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.04%  "

$ws.Range("D51").Value = "'0.4164"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.39%  "
